# The document contains four <id>...</id> tags (p038r_1 .. p038r_4), each
# split across three runs: "<id>", the bare id text, and "</id>". This
# collapses each trio of runs into a single run (keeping the formatting of
# the first "<id>" run) whose text is the full "<id>p038r_N</id>" string -
# mirroring a "newly downloaded tc" re-import that emits one run instead
# of three for this field.

$d = $word.ActiveDocument

$ids = @("p038r_1", "p038r_2", "p038r_3", "p038r_4")

foreach ($id in $ids) {
    $tag = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($tag, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $tag, 2) | Out-Null
}
